$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '76.871.31'
Set-TextValue 'E2' '  +0.76%  '
Set-TextValue 'D3' '2.957.72'
Set-TextValue 'E3' '  +2.79%  '
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '200.22'
Set-TextValue 'E5' '  +2.06%  '
Set-TextValue 'D6' '596.58'
Set-TextValue 'E6' '  -0.51%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  -0.13%  '
Set-TextValue 'E8' '  -0.83%  '
Set-TextValue 'E9' '  +1.55%  '
Set-TextValue 'D10' '2.956.00'
Set-TextValue 'E10' '  +2.91%  '
Set-TextValue 'E11' '  +13.97%  '
Set-TextValue 'E12' '  +0.43%  '
Set-TextValue 'D13' '3.500.93'
Set-TextValue 'E13' '  +2.74%  '
Set-TextValue 'E14' '  -0.30%  '
Set-TextValue 'D15' '76.739.93'
Set-TextValue 'E15' '  +0.69%  '
Set-TextValue 'D16' '28.28'
Set-TextValue 'E16' '  +2.77%  '
Set-TextValue 'E17' '  -0.76%  '
Set-TextValue 'D18' '2.957.11'
Set-TextValue 'E18' '  +2.68%  '
Set-TextValue 'D19' '13.46'
Set-TextValue 'E19' '  +7.60%  '
Set-TextValue 'D20' '8.71'
Set-TextValue 'E20' '  -3.65%  '
Set-TextValue 'D21' '372.95'
Set-TextValue 'E21' '  -2.69%  '
Set-TextValue 'D22' '4.32'
Set-TextValue 'E22' '  +4.37%  '
Set-TextValue 'D23' '2.26'
Set-TextValue 'E23' '  -3.30%  '
Set-TextValue 'D24' '72.68'
Set-TextValue 'E24' '  +1.28%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D25' '3.113.98'
Set-TextValue 'E25' '  +2.89%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'D27' '4.27'
Set-TextValue 'E27' '  +1.23%  '
Set-TextValue 'D28' '9.68'
Set-TextValue 'E28' '  -0.58%  '
Set-TextValue 'E29' '  +1.63%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  -0.19%  '
Set-TextValue 'D31' '8.20'
Set-TextValue 'E31' '  +5.98%  '
Set-TextValue 'D32' '1.38'
Set-TextValue 'E32' '  -1.70%  '
Set-TextValue 'D33' '496.20'
Set-TextValue 'E33' '  -3.33%  '
Set-TextValue 'E34' '  +1.04%  '
Set-TextValue 'E35' '  -0.19%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D36' '0.403'
Set-TextValue 'E36' '  +16.88%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D37' '166.22'
Set-TextValue 'E37' '  -0.55%  '
Set-TextValue 'D38' '0.113'
Set-TextValue 'E38' '  +23.38%  '
Set-TextValue 'D39' '20.18'
Set-TextValue 'E39' '  +0.64%  '
Set-TextValue 'E40' '  +1.34%  '
Set-TextValue 'E41' '  -6.69%  '
Set-TextValue 'E42' '  +0.04%  '
Set-TextValue 'D43' '180.92'
Set-TextValue 'E43' '  -2.41%  '
Set-TextValue 'D44' '4.92'
Set-TextValue 'E44' '  -2.90%  '
Set-TextValue 'E45' '  -1.53%  '
Set-TextValue 'D46' '40.13'
Set-TextValue 'E46' '  -0.35%  '
Set-TextValue 'E47' '  -4.07%  '
Set-TextValue 'E48' '  +2.29%  '
Set-TextValue 'D49' '3.88'
Set-TextValue 'E49' '  +3.51%  '
Set-TextValue 'E50' '  -2.64%  '
Set-TextValue 'D51' '22.62'
Set-TextValue 'E51' '  +4.89%  '
